$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# First summary block (Execute test suite) - rows 4-7, columns I:K
# ---------------------------------------------------------------------------
$ws.Range("I4:K4").Merge()
$ws.Range("I4").Value = "Result"
$ws.Range("I4:K4").Style = "60% - Accent2"
$ws.Range("I4:K4").HorizontalAlignment = -4108
$ws.Range("I4:K4").VerticalAlignment = -4108

$ws.Range("I5:J5").Merge()
$ws.Range("I5").Value = "สรุปผลการทดสอบ"
$ws.Range("I5:J5").Style = "40% - Accent2"
$ws.Range("I5:J5").HorizontalAlignment = -4108
$ws.Range("I5:J5").VerticalAlignment = -4108

$ws.Range("K5").Value = "คิดเป็น %"
$ws.Range("K5").Style = "40% - Accent2"
$ws.Range("K5").HorizontalAlignment = -4108
$ws.Range("K5").VerticalAlignment = -4108

$ws.Range("I6").Value = "Pass"
$ws.Range("J6").Formula = "=COUNTIF(F:F,""Pass"")"
$ws.Range("K6").Formula = "=TEXT(J6/3,""0.00%"")"
$ws.Range("I6:K6").Style = "Good"

$ws.Range("I7").Value = "Fail"
$ws.Range("J7").Formula = "=COUNTIF(F:F,""Fail"")"
$ws.Range("K7").Formula = "=TEXT(J7/3,""0.00%"")"
$ws.Range("I7:K7").Style = "Bad"

# ---------------------------------------------------------------------------
# Second summary block (Revise test suite) - rows 10-13, columns I:K
# ---------------------------------------------------------------------------
$ws.Range("I10:K10").Merge()
$ws.Range("I10").Value = "Revise"
$ws.Range("I10:K10").Style = "60% - Accent2"
$ws.Range("I10:K10").HorizontalAlignment = -4108
$ws.Range("I10:K10").VerticalAlignment = -4108

$ws.Range("I11:J11").Merge()
$ws.Range("I11").Value = "สรุปผลการทดสอบ"
$ws.Range("I11:J11").Style = "40% - Accent2"
$ws.Range("I11:J11").HorizontalAlignment = -4108
$ws.Range("I11:J11").VerticalAlignment = -4108

$ws.Range("K11").Value = "คิดเป็น %"
$ws.Range("K11").Style = "40% - Accent2"
$ws.Range("K11").HorizontalAlignment = -4108
$ws.Range("K11").VerticalAlignment = -4108

$ws.Range("I12").Value = "Pass"
$ws.Range("J12").Formula = "=COUNTIF(G:G,""Pass"")"
$ws.Range("K12").Formula = "=TEXT(J12/3,""0.00%"")"
$ws.Range("I12:K12").Style = "Good"

$ws.Range("I13").Value = "Fail"
$ws.Range("J13").Formula = "=COUNTIF(G:G,""Fail"")"
$ws.Range("K13").Formula = "=TEXT(J13/3,""0.00%"")"
$ws.Range("I13:K13").Style = "Bad"

# ---------------------------------------------------------------------------
# View state: reflect final selection after building the report
# ---------------------------------------------------------------------------
$ws.Range("K14").Select()
